# "output folder feature - extra context can be saved - repeated analyses are
#  saved separately": re-style the result table so that every data row (not
#  just the highlighted/filled ones) uses the Arial font with centered text.
#
# Rows 2,4,6,8,10,12 already carry a green fill (style index 1) - they keep
# that fill but gain the Arial font + centered alignment.
# Rows 3,5,7,9,11 have no fill - they just gain the Arial font + centered
# alignment.
#
# A scratch cell (far outside the used A1:I12 range) is used to build each
# combined format exactly once and then "paste special -> formats" it onto
# every row that needs it; this keeps the number of new style/font records
# written to styles.xml to a minimum (one combined format per case) instead
# of growing a new record per cell/row touched.

$wb = $wb
if (-not $wb) { $wb = $excel.ActiveWorkbook }
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlCenter = -4108

$fillRows  = @(2,4,6,8,10,12)
$plainRows = @(3,5,7,9,11)

# --- Combined format #1: Arial + centered + keep the existing green fill ---
$scratchFill = $ws.Range("K1")
$scratchFill.Interior.Color = $ws.Range("A2").Interior.Color
$scratchFill.Font.Name = "Arial"
$scratchFill.HorizontalAlignment = $xlCenter
$scratchFill.Copy()
foreach ($r in $fillRows) {
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial($xlPasteFormats)
}
$scratchFill.Clear()

# --- Combined format #2: Arial + centered, no fill ---
$scratchPlain = $ws.Range("K1")
$scratchPlain.Font.Name = "Arial"
$scratchPlain.HorizontalAlignment = $xlCenter
$scratchPlain.Copy()
foreach ($r in $plainRows) {
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial($xlPasteFormats)
}
$scratchPlain.Clear()
